$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" updates ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 854.3256832499998
$wsSchedule.Range("F2").Value = 14.12575534474206

# --- Sheet "Detailed" updates ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B19").Value = -5.97579
$wsDetailed.Range("B20").Value = -6.79084

$wsDetailed.Range("B21").Value = 0
$wsDetailed.Range("C21").Value = "historical"

$wsDetailed.Range("B22").Value = -2.78307
$wsDetailed.Range("C22").Value = "historical"

$wsDetailed.Range("B23").Value = -7.37741
$wsDetailed.Range("B24").Value = -8.19786
$wsDetailed.Range("B25").Value = -14
$wsDetailed.Range("B26").Value = -10
$wsDetailed.Range("B27").Value = -19.95
$wsDetailed.Range("B28").Value = -22.32715
$wsDetailed.Range("B29").Value = -10
$wsDetailed.Range("B30").Value = -22.12631
$wsDetailed.Range("B31").Value = -12.01
$wsDetailed.Range("B32").Value = -6.73554
$wsDetailed.Range("B33").Value = -7.38821
$wsDetailed.Range("B34").Value = -5.81015
$wsDetailed.Range("B35").Value = -4.61322

$wsDetailed.Range("B38").Value = 56.98

$wsDetailed.Range("B40").Value = 73.19
$wsDetailed.Range("B41").Value = 78
$wsDetailed.Range("B42").Value = 105
$wsDetailed.Range("B43").Value = 73.19
$wsDetailed.Range("B44").Value = 65
$wsDetailed.Range("B45").Value = 65

$wsDetailed.Range("B47").Value = 64.35272999999999

$wsDetailed.Range("B48").Value = 65
$wsDetailed.Range("B49").Value = 65
